$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATR")

# Rename the header of the first column (currently "index") to "i".
# This cell is the header of the "testdata" ListObject table, so updating
# the cell value also renames the table column.
$ws.Range("A1").Value = "i"

# The "index" column used to be 1-based (1..502 across rows 2..503);
# it is now 0-based (0..501).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}

# The "index" -> "i" header rename makes the best-fit column narrower
# (Excel recalculates the best-fit width to fit the shorter text).
$ws.Columns.Item(1).ColumnWidth = 3.15
